$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 648.14813
$ws.Range("J28").Value = 635
$ws.Range("L28").Value = 635
$ws.Range("N28").Value = -1605
$ws.Range("H43").Value = 10929.177
$ws.Range("I43").Value = 21541.334
$ws.Range("J43").Value = 5140.727
$ws.Range("K43").Value = 21541.334
$ws.Range("L43").Value = 5140.727
$ws.Range("M43").Value = -21472.334
$ws.Range("N43").Value = -5278.727
$ws.Range("H76").Value = 3245.3425
$ws.Range("I76").Value = 3131.6667
$ws.Range("J76").Value = 3770
$ws.Range("K76").Value = 3131.6667
$ws.Range("L76").Value = 3770
$ws.Range("M76").Value = -2816.6667
$ws.Range("N76").Value = -4400
$ws.Range("H79").Value = 3245.3425
$ws.Range("I79").Value = 3131.6667
$ws.Range("J79").Value = 3770
$ws.Range("K79").Value = 3131.6667
$ws.Range("L79").Value = 3770
$ws.Range("M79").Value = -2039.6667
$ws.Range("N79").Value = -5954
$ws.Range("H98").Value = 240.8
$ws.Range("I98").Value = 240.8
$ws.Range("K98").Value = 240.8
$ws.Range("M98").Value = 1257.2
$ws.Range("H106").Value = 27245.834
$ws.Range("I106").Value = 8983.333000000001
$ws.Range("J106").Value = 33333.332
$ws.Range("K106").Value = 8983.333000000001
$ws.Range("L106").Value = 33333.332
$ws.Range("M106").Value = -8352.333000000001
$ws.Range("N106").Value = -34595.332
$ws.Range("H122").Value = 240.8
$ws.Range("I122").Value = 240.8
$ws.Range("K122").Value = 722.4000000000001
$ws.Range("M122").Value = 1727.6
$ws.Range("H127").Value = 52965.1
$ws.Range("I127").Value = 52965.1
$ws.Range("K127").Value = 158895.3
$ws.Range("M127").Value = -153935.3
$ws.Range("H132").Value = 1325.6792
$ws.Range("I132").Value = 1016.22
$ws.Range("K132").Value = 3048.66
$ws.Range("M132").Value = -518.6599999999999
$ws.Range("H137").Value = 6713.8096
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 6713.8096
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("L137").Value = 20141.4288
$ws.Range("N137").Value = -25241.4288
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 6199.8335
$ws.Range("I45").Value = 5666.3335
$ws.Range("J45").Value = 6733.3335
$ws.Range("K45").Value = 5666.3335
$ws.Range("L45").Value = 6733.3335
$ws.Range("M45").Value = -5289.3335
$ws.Range("N45").Value = -7487.3335
$ws.Range("H110").Value = 1027.1936
$ws.Range("J110").Value = 1337.2222
$ws.Range("L110").Value = 1337.2222
$ws.Range("N110").Value = -5427.2222
$ws.Range("H132").Value = 1471893.6
$ws.Range("I132").Value = 1667822.8
$ws.Range("K132").Value = 5003468.4
$ws.Range("M132").Value = -5000938.4
$ws.Range("H134").Value = 81916.5
$ws.Range("J134").Value = 81916.5
$ws.Range("L134").Value = 81916.5
$ws.Range("N134").Value = -92056.5
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3488.2222
$ws.Range("I86").Value = 3182.5
$ws.Range("K86").Value = 3182.5
$ws.Range("M86").Value = -2059.5
$ws.Range("H89").Value = 3488.2222
$ws.Range("I89").Value = 3182.5
$ws.Range("K89").Value = 15912.5
$ws.Range("M89").Value = -10296.5
$ws.Range("H94").Value = 4027.0435
$ws.Range("I94").Value = 3525.1765
$ws.Range("J94").Value = 5449
$ws.Range("K94").Value = 3525.1765
$ws.Range("L94").Value = 5449
$ws.Range("M94").Value = -3074.1765
$ws.Range("N94").Value = -6351
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("N116").ClearContents()
$ws.Range("L116").Value = 0
$ws.Range("H118").Value = 58333.332
$ws.Range("J118").Value = 58333.332
$ws.Range("L118").Value = 58333.332
$ws.Range("N118").Value = -61647.332
$ws.Range("H132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("N132").ClearContents()
$ws.Range("L132").Value = 0
$ws.Range("H134").Value = 5598.971
$ws.Range("I134").Value = 2822.4707
$ws.Range("K134").Value = 8467.4121
$ws.Range("M134").Value = -5932.4121
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 364.85715
$ws.Range("I2").Value = 100
$ws.Range("J2").Value = 563.5
$ws.Range("K2").Value = 100
$ws.Range("L2").Value = 563.5
$ws.Range("M2").Value = 13
$ws.Range("N2").Value = -789.5
$ws.Range("H22").Value = 1948.2106
$ws.Range("I22").Value = 628
$ws.Range("K22").Value = 628
$ws.Range("M22").Value = -278
$ws.Range("H99").Value = 14146865
$ws.Range("I99").Value = 7429.5
$ws.Range("J99").Value = 19802640
$ws.Range("K99").Value = 7429.5
$ws.Range("L99").Value = 19802640
$ws.Range("M99").Value = -5931.5
$ws.Range("N99").Value = -19805636
$ws.Range("H108").Value = 60000
$ws.Range("J108").Value = 60000
$ws.Range("L108").Value = 60000
$ws.Range("N108").Value = -67680
$ws.Range("H126").Value = 14146865
$ws.Range("I126").Value = 7429.5
$ws.Range("J126").Value = 19802640
$ws.Range("K126").Value = 22288.5
$ws.Range("L126").Value = 59407920
$ws.Range("M126").Value = -19818.5
$ws.Range("N126").Value = -59412860
$ws.Range("H132").Value = 4366.6665
$ws.Range("I132").Value = 4000
$ws.Range("K132").Value = 12000
$ws.Range("M132").Value = -9470
$ws.Range("H134").Value = 3901.8333
$ws.Range("I134").Value = 3352.75
$ws.Range("K134").Value = 10058.25
$ws.Range("M134").Value = -7523.25
$ws.Range("H135").Value = 69999
$ws.Range("J135").Value = 69999
$ws.Range("L135").Value = 69999
$ws.Range("N135").Value = -80139
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 54.666668
$ws.Range("I2").Value = 46.714287
$ws.Range("J2").Value = 65.8
$ws.Range("K2").Value = 280.285722
$ws.Range("L2").Value = 394.8
$ws.Range("M2").Value = -167.285722
$ws.Range("N2").Value = -620.8
$ws.Range("H19").Value = 625.75
$ws.Range("J19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("N19").ClearContents()
$ws.Range("H22").Value = 1
$ws.Range("I22").Value = 1
$ws.Range("K22").Value = 3
$ws.Range("M22").Value = 166
$ws.Range("H27").Value = 1
$ws.Range("I27").Value = 1
$ws.Range("K27").Value = 3
$ws.Range("M27").Value = 99
$ws.Range("H103").Value = 327
$ws.Range("J103").Value = 361.2857
$ws.Range("L103").Value = 1083.8571
$ws.Range("N103").Value = -2841.8571
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H24").Value = 512499.5
$ws.Range("H29").Value = 3499.5
$ws.Range("J29").Value = 1999
$ws.Range("L29").Value = 1999
$ws.Range("N29").Value = -2579
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H126").Value = 14277.833
$ws.Range("I126").Value = 25139
$ws.Range("K126").Value = 75417
$ws.Range("M126").Value = -72947
$ws.Range("H132").Value = 11116.223
$ws.Range("I132").Value = 12863.571
$ws.Range("J132").Value = 5000.5
$ws.Range("K132").Value = 38590.713
$ws.Range("L132").Value = 15001.5
$ws.Range("M132").Value = -36060.713
$ws.Range("N132").Value = -20061.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3029.2222
$ws.Range("J22").Value = 3448.7334
$ws.Range("L22").Value = 3448.7334
$ws.Range("N22").Value = -4038.7334
$ws.Range("H27").Value = 3029.2222
$ws.Range("J27").Value = 3448.7334
$ws.Range("L27").Value = 3448.7334
$ws.Range("N27").Value = -3662.7334
$ws.Range("H40").Value = 2813.6538
$ws.Range("I40").Value = 2767.1052
$ws.Range("J40").Value = 2940
$ws.Range("K40").Value = 2767.1052
$ws.Range("L40").Value = 2940
$ws.Range("M40").Value = -2631.1052
$ws.Range("N40").Value = -3212
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("K63").Value = 0
$ws.Range("M63").ClearContents()
$ws.Range("H66").Value = 0
$ws.Range("I66").Value = 0
$ws.Range("K66").Value = 0
$ws.Range("M66").ClearContents()
$ws.Range("H119").Value = 99925
$ws.Range("J119").Value = 99925
$ws.Range("L119").Value = 99925
$ws.Range("N119").Value = -109601
$ws.Range("H132").Value = 3231.4
$ws.Range("I132").Value = 2892.5
$ws.Range("K132").Value = 8677.5
$ws.Range("M132").Value = -6147.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 23338334
$ws.Range("J11").Value = 25007500
$ws.Range("L11").Value = 25007500
$ws.Range("N11").Value = -25007784
$ws.Range("H62").Value = 6750
$ws.Range("J62").Value = 7500
$ws.Range("L62").Value = 7500
$ws.Range("N62").Value = -8748
$ws.Range("H65").Value = 6750
$ws.Range("J65").Value = 7500
$ws.Range("L65").Value = 37500
$ws.Range("N65").Value = -43740
$ws.Range("H126").Value = 2196.45
$ws.Range("I126").Value = 2008.125
$ws.Range("K126").Value = 6024.375
$ws.Range("M126").Value = -3554.375
$ws.Range("H132").Value = 2498.4285
$ws.Range("I132").Value = 2340.3076
$ws.Range("K132").Value = 7020.9228
$ws.Range("M132").Value = -4490.9228
